$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(5, 8).Value = 188.11765
$ws.Cells.Item(5, 9).Value = 137.92308
$ws.Cells.Item(5, 11).Value = 137.92308
$ws.Cells.Item(5, 13).Value = -22.92308

$ws.Cells.Item(92, 8).Value = 1709.2084
$ws.Cells.Item(92, 9).Value = 1632
$ws.Cells.Item(92, 11).Value = 1632
$ws.Cells.Item(92, 13).Value = -384

$ws.Cells.Item(93, 8).Value = 49995
$ws.Cells.Item(93, 10).Value = 49995
$ws.Cells.Item(93, 12).Value = 49995
$ws.Cells.Item(93, 14).Value = -54987

$ws.Cells.Item(98, 8).Value = 5273.724
$ws.Cells.Item(98, 9).Value = 4859.591
$ws.Cells.Item(98, 10).Value = 6575.2856
$ws.Cells.Item(98, 11).Value = 4859.591
$ws.Cells.Item(98, 12).Value = 6575.2856
$ws.Cells.Item(98, 13).Value = -3361.591
$ws.Cells.Item(98, 14).Value = -9571.285599999999

$ws.Cells.Item(112, 8).Value = 3295.0417
$ws.Cells.Item(112, 10).Value = 3575.2856
$ws.Cells.Item(112, 12).Value = 10725.8568
$ws.Cells.Item(112, 14).Value = -12941.8568

$ws.Cells.Item(122, 8).Value = 5273.724
$ws.Cells.Item(122, 9).Value = 4859.591
$ws.Cells.Item(122, 10).Value = 6575.2856
$ws.Cells.Item(122, 11).Value = 14578.773
$ws.Cells.Item(122, 12).Value = 19725.8568
$ws.Cells.Item(122, 13).Value = -12128.773
$ws.Cells.Item(122, 14).Value = -24625.8568

$ws.Cells.Item(132, 8).Value = 8065.488
$ws.Cells.Item(132, 9).Value = 6918.968
$ws.Cells.Item(132, 11).Value = 20756.904
$ws.Cells.Item(132, 13).Value = -18226.904

$ws.Cells.Item(133, 8).Value = 50155
$ws.Cells.Item(133, 10).Value = 50155
$ws.Cells.Item(133, 12).Value = 50155
$ws.Cells.Item(133, 14).Value = -60275

$ws.Cells.Item(135, 8).Value = 8108.4165
$ws.Cells.Item(135, 9).Value = 4263
$ws.Cells.Item(135, 10).Value = 15799.25
$ws.Cells.Item(135, 11).Value = 38367
$ws.Cells.Item(135, 12).Value = 142193.25
$ws.Cells.Item(135, 13).Value = -35832
$ws.Cells.Item(135, 14).Value = -147263.25

$ws.Cells.Item(140, 8).Value = 229992
$ws.Cells.Item(140, 10).Value = 229992
$ws.Cells.Item(140, 12).Value = 229992
$ws.Cells.Item(140, 14).Value = -240352

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 12105.241
$ws.Cells.Item(32, 9).Value = 10962.077
$ws.Cells.Item(32, 11).Value = 10962.077
$ws.Cells.Item(32, 13).Value = -10675.077

$ws.Cells.Item(74, 8).Value = 10464.348
$ws.Cells.Item(74, 9).Value = 1810.7028
$ws.Cells.Item(74, 11).Value = 1810.7028
$ws.Cells.Item(74, 13).Value = -936.7028

$ws.Cells.Item(77, 8).Value = 10464.348
$ws.Cells.Item(77, 9).Value = 1810.7028
$ws.Cells.Item(77, 11).Value = 9053.513999999999
$ws.Cells.Item(77, 13).Value = -4685.513999999999

$ws.Cells.Item(130, 8).Value = 51159.25
$ws.Cells.Item(130, 10).Value = 51159.25
$ws.Cells.Item(130, 12).Value = 51159.25
$ws.Cells.Item(130, 14).Value = -61199.25

$ws.Cells.Item(132, 8).Value = 11065.182
$ws.Cells.Item(132, 9).Value = 1350.6765
$ws.Cells.Item(132, 11).Value = 4052.0295
$ws.Cells.Item(132, 13).Value = -1522.0295

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 62502744
$ws.Cells.Item(86, 9).Value = 2800.375
$ws.Cells.Item(86, 10).Value = 125002690
$ws.Cells.Item(86, 11).Value = 2800.375
$ws.Cells.Item(86, 12).Value = 125002690
$ws.Cells.Item(86, 13).Value = -1677.375
$ws.Cells.Item(86, 14).Value = -125004936

$ws.Cells.Item(89, 8).Value = 62502744
$ws.Cells.Item(89, 9).Value = 2800.375
$ws.Cells.Item(89, 10).Value = 125002690
$ws.Cells.Item(89, 11).Value = 14001.875
$ws.Cells.Item(89, 12).Value = 625013450
$ws.Cells.Item(89, 13).Value = -8385.875
$ws.Cells.Item(89, 14).Value = -625024682

$ws.Cells.Item(94, 8).Value = 3100.8215
$ws.Cells.Item(94, 9).Value = 2051.0454
$ws.Cells.Item(94, 11).Value = 2051.0454
$ws.Cells.Item(94, 13).Value = -1600.0454

$ws.Cells.Item(126, 8).Value = 85000
$ws.Cells.Item(126, 10).Value = 85000
$ws.Cells.Item(126, 12).Value = 85000
$ws.Cells.Item(126, 14).Value = -94880

$ws.Cells.Item(132, 8).Value = 78192.25
$ws.Cells.Item(132, 10).Value = 78192.25
$ws.Cells.Item(132, 12).Value = 78192.25
$ws.Cells.Item(132, 14).Value = -88312.25

$ws.Cells.Item(134, 8).Value = 67052.39999999999
$ws.Cells.Item(134, 9).Value = 92360
$ws.Cells.Item(134, 10).Value = 29091
$ws.Cells.Item(134, 11).Value = 277080
$ws.Cells.Item(134, 12).Value = 87273
$ws.Cells.Item(134, 13).Value = -274545
$ws.Cells.Item(134, 14).Value = -92343

$ws.Cells.Item(135, 8).Value = 73130
$ws.Cells.Item(135, 10).Value = 73130
$ws.Cells.Item(135, 12).Value = 73130
$ws.Cells.Item(135, 14).Value = -83270

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(4, 8).Value = 353999.7
$ws.Cells.Item(4, 10).Value = 999999
$ws.Cells.Item(4, 12).Value = 999999
$ws.Cells.Item(4, 14).Value = -1000223

$ws.Cells.Item(31, 8).Value = 56394.34
$ws.Cells.Item(31, 9).Value = 92719.69500000001
$ws.Cells.Item(31, 10).Value = 16609.428
$ws.Cells.Item(31, 11).Value = 92719.69500000001
$ws.Cells.Item(31, 12).Value = 16609.428
$ws.Cells.Item(31, 13).Value = -92424.69500000001
$ws.Cells.Item(31, 14).Value = -17199.428

$ws.Cells.Item(34, 8).Value = 56394.34
$ws.Cells.Item(34, 9).Value = 92719.69500000001
$ws.Cells.Item(34, 10).Value = 16609.428
$ws.Cells.Item(34, 11).Value = 92719.69500000001
$ws.Cells.Item(34, 12).Value = 16609.428
$ws.Cells.Item(34, 13).Value = -92517.69500000001
$ws.Cells.Item(34, 14).Value = -17013.428

$ws.Cells.Item(69, 8).Value = 16250
$ws.Cells.Item(69, 9).Value = 16250
$ws.Cells.Item(69, 11).Value = 16250
$ws.Cells.Item(69, 13).Value = -15501

$ws.Cells.Item(72, 8).Value = 16250
$ws.Cells.Item(72, 9).Value = 16250
$ws.Cells.Item(72, 11).Value = 48750
$ws.Cells.Item(72, 13).Value = -45006

$ws.Cells.Item(109, 8).Value = 14583.333
$ws.Cells.Item(109, 10).Value = 14583.333
$ws.Cells.Item(109, 12).Value = 14583.333
$ws.Cells.Item(109, 14).Value = -16663.333

$ws.Cells.Item(122, 8).Value = 3648.1177
$ws.Cells.Item(122, 9).Value = 1029.875
$ws.Cells.Item(122, 10).Value = 5975.4443
$ws.Cells.Item(122, 11).Value = 3089.625
$ws.Cells.Item(122, 12).Value = 17926.3329
$ws.Cells.Item(122, 13).Value = -639.625
$ws.Cells.Item(122, 14).Value = -22826.3329

$ws.Cells.Item(132, 8).Value = 4607.2793
$ws.Cells.Item(132, 9).Value = 1533.3055
$ws.Cells.Item(132, 11).Value = 4599.916499999999
$ws.Cells.Item(132, 13).Value = -2069.916499999999

$ws.Cells.Item(134, 8).Value = 5419.073
$ws.Cells.Item(134, 9).Value = 2219.6333
$ws.Cells.Item(134, 11).Value = 6658.8999
$ws.Cells.Item(134, 13).Value = -4123.8999

$ws.Cells.Item(139, 8).Value = 124500
$ws.Cells.Item(139, 10).Value = 124500
$ws.Cells.Item(139, 12).Value = 124500
$ws.Cells.Item(139, 14).Value = -134780

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 3477420.5
$ws.Cells.Item(4, 10).Value = 1000000
$ws.Cells.Item(4, 12).Value = 3000000
$ws.Cells.Item(4, 14).Value = -3000224

$ws.Cells.Item(34, 8).Value = 3967
$ws.Cells.Item(34, 10).Value = 6032
$ws.Cells.Item(34, 12).Value = 18096
$ws.Cells.Item(34, 14).Value = -18264

$ws.Cells.Item(39, 8).Value = 5537.375
$ws.Cells.Item(39, 10).Value = 5566.6665
$ws.Cells.Item(39, 12).Value = 16699.9995
$ws.Cells.Item(39, 14).Value = -17287.9995

$ws.Cells.Item(46, 8).Value = 330.92307
$ws.Cells.Item(46, 9).Value = 330.92307
$ws.Cells.Item(46, 10).Value = 0
$ws.Cells.Item(46, 11).Value = 992.7692099999999
$ws.Cells.Item(46, 12).Value = 0
$ws.Cells.Item(46, 13).Value = -901.7692099999999
$ws.Cells.Item(46, 14).ClearContents() | Out-Null

$ws.Cells.Item(107, 8).Value = 842.2414
$ws.Cells.Item(107, 10).Value = 995.8823
$ws.Cells.Item(107, 12).Value = 2987.6469
$ws.Cells.Item(107, 14).Value = -6827.6469

$ws.Cells.Item(129, 8).Value = 4786788.5
$ws.Cells.Item(129, 9).Value = 1789
$ws.Cells.Item(129, 11).Value = 5367
$ws.Cells.Item(129, 13).Value = -367

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(104, 8).Value = 45917.75
$ws.Cells.Item(104, 10).Value = 45917.75
$ws.Cells.Item(104, 12).Value = 45917.75
$ws.Cells.Item(104, 14).Value = -52905.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(130, 8).Value = 60000
$ws.Cells.Item(130, 10).Value = 60000
$ws.Cells.Item(130, 12).Value = 60000
$ws.Cells.Item(130, 14).Value = -70040

$ws.Cells.Item(136, 8).Value = 37050.242
$ws.Cells.Item(136, 9).Value = 40273.367
$ws.Cells.Item(136, 10).Value = 20290
$ws.Cells.Item(136, 11).Value = 120820.101
$ws.Cells.Item(136, 12).Value = 60870
$ws.Cells.Item(136, 13).Value = -118270.101
$ws.Cells.Item(136, 14).Value = -65970

$ws.Cells.Item(139, 8).Value = 88205.45
$ws.Cells.Item(139, 9).Value = 37575
$ws.Cells.Item(139, 11).Value = 37575
$ws.Cells.Item(139, 13).Value = -32435

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 120349.91
$ws.Cells.Item(62, 9).Value = 18466.5
$ws.Cells.Item(62, 11).Value = 18466.5
$ws.Cells.Item(62, 13).Value = -17842.5

$ws.Cells.Item(65, 8).Value = 120349.91
$ws.Cells.Item(65, 9).Value = 18466.5
$ws.Cells.Item(65, 11).Value = 92332.5
$ws.Cells.Item(65, 13).Value = -89212.5

$ws.Cells.Item(81, 8).Value = 965.3333
$ws.Cells.Item(81, 9).Value = 965.3333
$ws.Cells.Item(81, 10).Value = 0
$ws.Cells.Item(81, 11).Value = 1930.6666
$ws.Cells.Item(81, 12).Value = 0
$ws.Cells.Item(81, 13).Value = -869.6666
$ws.Cells.Item(81, 14).ClearContents() | Out-Null

$ws.Cells.Item(84, 8).Value = 965.3333
$ws.Cells.Item(84, 9).Value = 965.3333
$ws.Cells.Item(84, 10).Value = 0
$ws.Cells.Item(84, 11).Value = 9653.333000000001
$ws.Cells.Item(84, 12).Value = 0
$ws.Cells.Item(84, 13).Value = -4349.333000000001
$ws.Cells.Item(84, 14).ClearContents() | Out-Null

$ws.Cells.Item(122, 8).Value = 589114.4399999999
$ws.Cells.Item(122, 9).Value = 856032.9399999999
$ws.Cells.Item(122, 11).Value = 2568098.82
$ws.Cells.Item(122, 13).Value = -2565648.82

$ws.Cells.Item(130, 8).Value = 45000
$ws.Cells.Item(130, 10).Value = 45000
$ws.Cells.Item(130, 12).Value = 45000
$ws.Cells.Item(130, 14).Value = -55040

$ws.Cells.Item(132, 8).Value = 7384.489
$ws.Cells.Item(132, 9).Value = 2396.9565
$ws.Cells.Item(132, 10).Value = 12598.728
$ws.Cells.Item(132, 11).Value = 7190.869499999999
$ws.Cells.Item(132, 12).Value = 37796.18399999999
$ws.Cells.Item(132, 13).Value = -4660.869499999999
$ws.Cells.Item(132, 14).Value = -42856.18399999999

$ws.Cells.Item(141, 8).Value = 65242
$ws.Cells.Item(141, 10).Value = 65242
$ws.Cells.Item(141, 12).Value = 65242
$ws.Cells.Item(141, 14).Value = -75602
